# Generate Report for Handoff
#
# The f337a311-... source file has finished its localization round trip and
# drops out of the report entirely (its row is removed from every sheet).
# The remaining 18fa94ca-... row moves back into the handoff queue: its
# Status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" and its Latest Handoff Datetime is refreshed.

$wb = $excel.ActiveWorkbook

$localizationConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b4b66ae46207e90157b73158cd37fa5f24b8e782/.localization-config"

function Delete-HyperlinksByDisplay($sheet, $displayText) {
    $guard = 0
    $continueLoop = $true
    while ($continueLoop -and $guard -lt 50) {
        $guard++
        $found = $null
        foreach ($h in $sheet.Hyperlinks) {
            if ($h.TextToDisplay -eq $displayText) {
                $found = $h
                break
            }
        }
        if ($found -ne $null) {
            $found.Delete()
        } else {
            $continueLoop = $false
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

Delete-HyperlinksByDisplay $overview "f337a311-5fd0-4ad7-b2ea-4019cc13ed0d.md"
Delete-HyperlinksByDisplay $overview ".localization-config"

$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Rows.Item(3).Delete()

$overview.Hyperlinks.Add($overview.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

Delete-HyperlinksByDisplay $zhcn "f337a311-5fd0-4ad7-b2ea-4019cc13ed0d.md"
Delete-HyperlinksByDisplay $zhcn "f337a311-5fd0-4ad7-b2ea-4019cc13ed0d.c95a77cb2ca4cd6f3be82e6d6874e59cd2be5b5d.zh-cn.xlf"
Delete-HyperlinksByDisplay $zhcn ".localization-config"

$zhcn.Range("B2").Value = "Ready for handoff"
$zhcn.Range("D2").Value = "2016-02-22 04:43:38"
$zhcn.Rows.Item(3).Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

Delete-HyperlinksByDisplay $dede "f337a311-5fd0-4ad7-b2ea-4019cc13ed0d.md"
Delete-HyperlinksByDisplay $dede "f337a311-5fd0-4ad7-b2ea-4019cc13ed0d.c95a77cb2ca4cd6f3be82e6d6874e59cd2be5b5d.de-de.xlf"
Delete-HyperlinksByDisplay $dede ".localization-config"

$dede.Range("B2").Value = "Ready for handoff"
$dede.Range("D2").Value = "2016-02-22 04:43:52"
$dede.Rows.Item(3).Delete()

$dede.Hyperlinks.Add($dede.Range("A3"), $localizationConfigUrl, "", "", ".localization-config") | Out-Null
